$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 21 corresponds to month 45870 (2025-08) - update stats
$ws.Range("B21").Value = 6219
$ws.Range("C21").Value = 980
$ws.Range("D21").Value = 5575337
$ws.Range("E21").Value = 896.5005627914455
$ws.Range("F21").Value = 7.950008679048781
$ws.Range("G21").Value = 3.375527426160341
$ws.Range("H21").Value = 27.2625572589827
